# "deck select on job error more obvious"
# Add a new error-code row to the Error sheet's table:
#   1004 -> "一套卡牌不能有两种职业卡牌" (a deck cannot contain two different class' cards)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row right after the existing last row (row 14).
$ws.Range("A15").Value = 1004
$ws.Range("B15").Value = "一套卡牌不能有两种职业卡牌"

# Grow the worksheet's table ("表1") so the new row is included in it,
# which keeps the table ref / autoFilter ranges in sync (A1:B14 -> A1:B15).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B15")) | Out-Null

# Match Excel's natural behavior of leaving the selection on the last
# edited cell (B15) after the new row is entered.
$ws.Range("B15").Select() | Out-Null
